$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1996.6666
$ws.Range("J17").Value = 1996.6666
$ws.Range("L17").Value = 5989.9998
$ws.Range("N17").Value = -6325.9998

$ws.Range("H28").Value = 1738
$ws.Range("I28").Value = 1299
$ws.Range("K28").Value = 1299
$ws.Range("M28").Value = -814

$ws.Range("H74").Value = 4599.8
$ws.Range("I74").Value = 4599.8
$ws.Range("K74").Value = 4599.8
$ws.Range("M74").Value = -3663.8

$ws.Range("H77").Value = 4599.8
$ws.Range("I77").Value = 4599.8
$ws.Range("K77").Value = 22999
$ws.Range("M77").Value = -18319

$ws.Range("H138").Value = 3300.0686
$ws.Range("I138").Value = 2099.4285
$ws.Range("J138").Value = 3427.4092
$ws.Range("K138").Value = 6298.2855
$ws.Range("L138").Value = 10282.2276
$ws.Range("M138").Value = -1158.2855
$ws.Range("N138").Value = -20562.2276

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10979.805
$ws.Range("I32").Value = 9737.5625
$ws.Range("K32").Value = 9737.5625
$ws.Range("M32").Value = -9450.5625

$ws.Range("H102").Value = 1459.1666
$ws.Range("I102").Value = 1699.8
$ws.Range("K102").Value = 1699.8
$ws.Range("M102").Value = -77.79999999999995

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 345.875
$ws.Range("I22").Value = 178.6
$ws.Range("K22").Value = 178.6
$ws.Range("M22").Value = 171.4

$ws.Range("H31").Value = 4669.129
$ws.Range("I31").Value = 3113.2632
$ws.Range("J31").Value = 7132.5835
$ws.Range("K31").Value = 3113.2632
$ws.Range("L31").Value = 7132.5835
$ws.Range("M31").Value = -2818.2632
$ws.Range("N31").Value = -7722.5835

$ws.Range("H34").Value = 4669.129
$ws.Range("I34").Value = 3113.2632
$ws.Range("J34").Value = 7132.5835
$ws.Range("K34").Value = 3113.2632
$ws.Range("L34").Value = 7132.5835
$ws.Range("M34").Value = -2911.2632
$ws.Range("N34").Value = -7536.5835

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 563.3125
$ws.Range("J5").Value = 556.1818
$ws.Range("L5").Value = 1668.5454
$ws.Range("N5").Value = -1892.5454

$ws.Range("H62").Value = 4600
$ws.Range("I62").Value = 8500
$ws.Range("K62").Value = 25500
$ws.Range("M62").Value = -24814

$ws.Range("H64").Value = 0
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = $null
$ws.Range("M64").Value = $null
$ws.Range("N64").Value = 0

$ws.Range("H65").Value = 4600
$ws.Range("I65").Value = 8500
$ws.Range("K65").Value = 76500
$ws.Range("M65").Value = -73068

$ws.Range("H67").Value = 0
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = $null
$ws.Range("M67").Value = $null
$ws.Range("N67").Value = 0

$ws.Range("H68").Value = 1692
$ws.Range("I68").Value = 1447.8
$ws.Range("J68").Value = 1997.25
$ws.Range("K68").Value = 4343.4
$ws.Range("L68").Value = 5991.75
$ws.Range("M68").Value = -3532.4
$ws.Range("N68").Value = -7613.75

$ws.Range("H69").Value = 1887.6666
$ws.Range("J69").Value = 1498.625
$ws.Range("L69").Value = 4495.875
$ws.Range("N69").Value = -6117.875

$ws.Range("H70").Value = 4000
$ws.Range("I70").Value = 4000
$ws.Range("K70").Value = 12000
$ws.Range("M70").Value = -11685

$ws.Range("H71").Value = 1692
$ws.Range("I71").Value = 1447.8
$ws.Range("J71").Value = 1997.25
$ws.Range("K71").Value = 13030.2
$ws.Range("L71").Value = 17975.25
$ws.Range("M71").Value = -8974.199999999999
$ws.Range("N71").Value = -26087.25

$ws.Range("H72").Value = 1887.6666
$ws.Range("J72").Value = 1498.625
$ws.Range("L72").Value = 13487.625
$ws.Range("N72").Value = -21599.625

$ws.Range("H73").Value = 4000
$ws.Range("I73").Value = 4000
$ws.Range("K73").Value = 12000
$ws.Range("M73").Value = -10908

$ws.Range("H81").Value = 2996.5715
$ws.Range("J81").Value = 2996.5715
$ws.Range("L81").Value = 8989.7145
$ws.Range("N81").Value = -11235.7145

$ws.Range("H82").Value = 16666.666
$ws.Range("J82").Value = 20000
$ws.Range("L82").Value = 60000
$ws.Range("N82").Value = -60812

$ws.Range("H84").Value = 2996.5715
$ws.Range("J84").Value = 2996.5715
$ws.Range("L84").Value = 26969.1435
$ws.Range("N84").Value = -38201.1435

$ws.Range("H85").Value = 16666.666
$ws.Range("J85").Value = 20000
$ws.Range("L85").Value = 60000
$ws.Range("N85").Value = -62808

$ws.Range("H102").Value = 6166.3335
$ws.Range("I102").Value = 6166.3335
$ws.Range("K102").Value = 18499.0005
$ws.Range("M102").Value = -16065.0005

$ws.Range("H105").Value = 7920
$ws.Range("J105").Value = 7920
$ws.Range("L105").Value = 23760
$ws.Range("N105").Value = -29002

$ws.Range("H106").Value = 13712.714
$ws.Range("J106").Value = 13712.714
$ws.Range("L106").Value = 41138.142
$ws.Range("N106").Value = -43030.142

$ws.Range("H135").Value = 563.3125
$ws.Range("J135").Value = 556.1818
$ws.Range("L135").Value = 5005.6362
$ws.Range("N135").Value = -10075.6362

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2054.8333
$ws.Range("I102").Value = 582.3333
$ws.Range("K102").Value = 582.3333
$ws.Range("M102").Value = 1039.6667

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 5010.357
$ws.Range("I81").Value = 4023.5454
$ws.Range("J81").Value = 8628.666999999999
$ws.Range("K81").Value = 8047.0908
$ws.Range("L81").Value = 17257.334
$ws.Range("M81").Value = -6986.0908
$ws.Range("N81").Value = -19379.334

$ws.Range("H84").Value = 5010.357
$ws.Range("I84").Value = 4023.5454
$ws.Range("J84").Value = 8628.666999999999
$ws.Range("K84").Value = 40235.454
$ws.Range("L84").Value = 86286.67
$ws.Range("M84").Value = -34931.454
$ws.Range("N84").Value = -96894.67

$ws.Range("H125").Value = 149125
$ws.Range("J125").Value = 149125
$ws.Range("L125").Value = 149125
$ws.Range("N125").Value = -158965
